$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Meeting-info table (Table 1): tighten column widths by 1 twip
#    (2233 -> 2232, 2030 -> 2031; the gridSpan=4 header cell's own width
#    cannot be addressed independently of column 2 in this engine, so it is
#    left to whatever the column-2 resize produces).
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Cell(2, 1).Width = 111.6   # 2232 dxa
$t1.Cell(2, 2).Width = 101.55  # 2031 dxa

# Date of meeting: "../09/2025" -> "21/09/2025"
$d.Content.Find.Execute("../09/2025", $false, $false, $false, $false, $false, `
                         $true, 0, $false, "21/09/2025", 1) | Out-Null

# ---------------------------------------------------------------------------
# 2) "...separated it into 6 main functions..." - the " " run and the
#    "6 main " run become a single " 6 main " run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" 6 main ", $false, $false, $false, $false, $false, `
                         $true, 0, $false, " 6 main ", 1) | Out-Null

# ---------------------------------------------------------------------------
# 3) "3. ACTION ITEMS" table (Table 4): tighten column widths by 1 twip.
# ---------------------------------------------------------------------------
$t4 = $d.Tables.Item(4)
$t4.Cell(2, 1).Width = 103.65  # 2073 dxa
$t4.Cell(2, 2).Width = 122.55  # 2451 dxa

# ---------------------------------------------------------------------------
# 4) Second action-items table (Table 5): tighten column widths by 1 twip
#    and collapse the split "2"+"3"+"/09/2025" runs into "23/09/2025" for
#    the first two data rows only.
# ---------------------------------------------------------------------------
$t5 = $d.Tables.Item(5)
$t5.Cell(1, 1).Width = 103.65  # 2073 dxa
$t5.Cell(1, 2).Width = 122.55  # 2451 dxa

$c1 = $t5.Cell(1, 3).Range
$c1.Find.Execute("23/09/2025", $false, $false, $false, $false, $false, `
                  $true, 0, $false, "23/09/2025", 1) | Out-Null

$c2 = $t5.Cell(2, 3).Range
$c2.Find.Execute("23/09/2025", $false, $false, $false, $false, $false, `
                  $true, 0, $false, "23/09/2025", 1) | Out-Null

Write-Output "done"
